$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "release/8.0.2"
$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("D5").Value = "X"
$ws.Range("E5").Value = "X"
